$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.256.92"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "3.517.63"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "610.66"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").Value = "148.32"
$ws.Range("E6").Value = "  -1.63%  "

$ws.Range("D7").Value = "3.516.36"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -1.53%  "

$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").Value = "8.03"
$ws.Range("E11").Value = "  +6.44%  "

$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").Value = "0.0000217"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").Value = "4.112.25"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "31.55"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "3.523.55"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "67.240.46"
$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "10.95"
$ws.Range("E19").Value = "  +9.47%  "

$ws.Range("E20").Value = "  -2.21%  "

$ws.Range("D21").Value = "15.44"
$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").Value = "437.07"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").Value = "0.610"
$ws.Range("E23").Value = "  -2.85%  "

$ws.Range("D24").Value = "80.09"
$ws.Range("E24").Value = "  +1.37%  "

$ws.Range("D25").Value = "3.657.38"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -3.90%  "

$ws.Range("D28").Value = "9.84"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("E29").Value = "  -5.04%  "

$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("E31").Value = "  -4.59%  "

$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -1.79%  "

$ws.Range("D34").Value = "25.65"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "5.98"
$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").Value = "8.06"
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D40").Value = "176.30"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("D41").Value = "0.0900"

$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "2.05"
$ws.Range("E43").Value = "  -10.20%  "

$ws.Range("D44").Value = "0.898"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "46.40"
$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("D46").Value = "28.10"
$ws.Range("E46").Value = "  -8.68%  "

$ws.Range("E47").Value = "  -4.56%  "

$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("E49").Value = "  -1.69%  "

$ws.Range("D50").Value = "0.996"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("E51").Value = "  -2.08%  "
